# Auto-generated Excel COM-interop edit script
# Applies grading updates ("complementos posteriores a la entrega") to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Feedback comment cells (column G): updated grading commentary ---
$ws.Range("G192").Value = "Comentario personal `"función`":`n- Aunque no exploras mucho en posibilidades para relacionarte con el entorno (eso podría estar en muchos lugares), lo que estás proponiendo funciona.`n- Al explorar poco por lo que indicamos en el comentario anterior, es necesario que te exijas más al momento de desarrollar los diferentes dibujos abordando más detalles y complejidad de estos.`n- Te hace falta incorporar estructura e indicar ejes en planta.`n- La historia de tu familia es un insumo esencial para validar lo que estás proponiendo para esa familia. Sin embargo, no la has cargado a la carpeta de TEAMS.`n`nComentario personal `"representación`":`n- Te hace falta incorporar la planta en 1:50. Esta escala te permitirá resolver detalles que solo se ven en ese tamaño.`n- El espacio que estás considerando para tu sección por fachada es insuficiente para lo que necesitas.`n- El 3D q presentas no corresponde a lo que estás diseñando. Esto es grave!`n- Debes tener intenciones para diagramar más claras y efectivas. Esto lo podremos revisar en una próxima asesoría, pero por lo pronto, inclinar de manera aleatoria dibujos cuya representación técnica surge de una vista ortogonal (ej: planta, sección y fachada) no le aporta valor a tu memoria y antes confunde.`n`nLa nota general para esta entrega es 3,8."
$ws.Range("G217").Value = "Comentario general equipo:`n- Es perfectamente factible para esta entrega, haber entregado hojas separadas; es importante sin embargo, que activen ya mismo, la diagramación de esos dibujos en clave de todo el conjunto. en este momento que empiecen a pensar en clave de diagramación de todo el conjunto.`n`nComentario personal función:`n- Te hace falta cuadro de áreas`n- A nivel de planta del primer nivel está muy afinado el desarrollo de tu vivienda.`n- Es necesario que te hagas preguntas desde lo que sucede en alzado (fachada y sección): cómo enriqueces lo que tienes en planta con lo que sucede en la cubierta?`n- Te hace falta indicar pendientes en planta de cubierta y nombres de espacios en el primer nivel.`n- La historia de tu familia es un insumo esencial para validar lo que estás proponiendo para esa familia. Sin embargo, no la has cargado a la carpeta de TEAMS.`n`nComentario personal representación:`n- Debes considerar las escalas de los planos solicitados: entregas una planta en 1:200 y necesitabas presentarla en 1:50. Al no representar en dichas escalas, no podemos ver el detalle de ciertos aspectos del dibujo que solo se ven a ese tamaño.`n- Le das mucha más importancia al 3D en el espacio total disponible para la memoria, pero en contra del espacio que requieres para la planta en un nivel.`n`nLa nota general para esta entrega es 4,4."
$ws.Range("G642").Value = "Al revisar lo que entregaste en una fecha posterior, debido al apoyo que le diste de manera oportuna a tu compañera, hemos ajustado tu evaluación. Aún es necesario que eleves tu nivel de reflexión espacial y producción para los productos que debes presentar en la entrega final.`nLa nota general para esta entrega es 4.3."
$ws.Range("G842").Value = "Comentario personal `"función`":`n- La resolución que tienen los dibujos que presentas, dificulta mucho poder revisar detalles de la funcionalidad de tu diseño: debemos verlo en detalle en una próxima asesoría.`n- Las intenciones que tienes con los dos volúmenes de tu vivienda, podrían dialogar mejor entre sí: tener intenciones a nivel de conjunto (en este momento son dos volúmenes que no dialogan formalmente).`n- La historia de tu familia es un insumo esencial para validar lo que estás proponiendo para esa familia. Sin embargo, no la has cargado a la carpeta de TEAMS.`n`nComentario personal `"representación`":`n- Revisa intenciones de elementos de diagramación pues generan conflictos para leer planimetría con claridad (ej: fondo azul.`n- Revisa escala de planos que presentas: para la próxima semana, vuelve a traer impresas 3 plantas en Esc-1:50 para que podamos resolver detalles funcionales (no incluyas detalles de ambientación de la memoria: solo planos).`n- Revisa intenciones al representar y la cantidad de cosas que añades.`n- Aumenta el tamaño del cuadro de áreas para poder ver su información.`nComentario personal `"argumentación`":`n- Has subido la historia de tu familia posterior a la entrega, pero lo has hecho a manera de fotografía desde tu bitácota y les habíamos solicitado que digitalizaran dicho texto (ej: word)`nNota general:`nTu nota general es 3.9 y por tener 0.5 acumulado de la entrega previa, esta sube a 4,4"
$ws.Range("G867").Value = "Al revisar lo que entregaste en una fecha posterior, debido a la situación difícil que viviste, hemos ajustado tu evaluación. Reconocemos tu valor para continuar y la capacidad para reintegrarte a las actividades del taller. Aún es necesario que eleves tu nivel de reflexión espacial y producción para los productos que debes presentar en la entrega final.`nLa nota general para esta entrega es 4.3"

# --- Rubric rows: Nivel / Nivel_eq / Categoria_match / Descrip ---
# Nivel (column C) must stay a TEXT value (e.g. "5.0"), not a number,
# so force Text number format before assigning, then restore the default
# style afterward so no stray number-format style is left on the cell.
$ws.Range("C195").NumberFormat = "@"
$ws.Range("C195").Value = "5.0"
$ws.Range("C195").Style = "Normal"
$ws.Range("D195").Value = "A"
$ws.Range("F195").Value = "E4.Arg_A"
$ws.Range("G195").Value = "Argumenta las decisiones de diseño de manera sólida y coherente, basándose en fundamentos disciplinares y un profundo conocimiento del usuario. Las decisiones están bien justificadas y son consistentes."

$ws.Range("C220").NumberFormat = "@"
$ws.Range("C220").Value = "5.0"
$ws.Range("C220").Style = "Normal"
$ws.Range("D220").Value = "A"
$ws.Range("F220").Value = "E4.Arg_A"
$ws.Range("G220").Value = "Argumenta las decisiones de diseño de manera sólida y coherente, basándose en fundamentos disciplinares y un profundo conocimiento del usuario. Las decisiones están bien justificadas y son consistentes."

$ws.Range("C643").NumberFormat = "@"
$ws.Range("C643").Value = "3.5"
$ws.Range("C643").Style = "Normal"
$ws.Range("D643").Value = "B"
$ws.Range("F643").Value = "E4.Func_B"
$ws.Range("G643").Value = "Aplica estrategias de composición formal en una vivienda con algunas inconsistencias en la integración de las variables propuestas para el ejercicio. El proyecto presenta algunas debilidades en su función y/o en su resolución técnica."

$ws.Range("C644").NumberFormat = "@"
$ws.Range("C644").Value = "3.8"
$ws.Range("C644").Style = "Normal"
$ws.Range("D644").Value = "B"
$ws.Range("F644").Value = "E4.Rep_B"
$ws.Range("G644").Value = "Utiliza recursos verbales y gráficos en 2 y 3 dimensiones para representar y comunicar el proyecto, aunque con algunas inconsistencias en la coherencia y la efectividad de la comunicación. La presentación es generalmente clara, pero ciertos aspectos claves del proyecto no están suficientemente desarrollados o explicados."

$ws.Range("C845").NumberFormat = "@"
$ws.Range("C845").Value = "3.0"
$ws.Range("C845").Style = "Normal"
$ws.Range("D845").Value = "B"
$ws.Range("F845").Value = "E4.Arg_B"
$ws.Range("G845").Value = "Presenta una argumentación aceptable para las decisiones de diseño, pero con algunas inconsistencias o falta de profundidad en los fundamentos disciplinares o el conocimiento del usuario."

$ws.Range("C868").NumberFormat = "@"
$ws.Range("C868").Value = "3.5"
$ws.Range("C868").Style = "Normal"
$ws.Range("D868").Value = "B"
$ws.Range("F868").Value = "E4.Func_B"
$ws.Range("G868").Value = "Aplica estrategias de composición formal en una vivienda con algunas inconsistencias en la integración de las variables propuestas para el ejercicio. El proyecto presenta algunas debilidades en su función y/o en su resolución técnica."

$ws.Range("C869").NumberFormat = "@"
$ws.Range("C869").Value = "3.8"
$ws.Range("C869").Style = "Normal"
$ws.Range("D869").Value = "B"
$ws.Range("F869").Value = "E4.Rep_B"
$ws.Range("G869").Value = "Utiliza recursos verbales y gráficos en 2 y 3 dimensiones para representar y comunicar el proyecto, aunque con algunas inconsistencias en la coherencia y la efectividad de la comunicación. La presentación es generalmente clara, pero ciertos aspectos claves del proyecto no están suficientemente desarrollados o explicados."

$ws.Range("C870").NumberFormat = "@"
$ws.Range("C870").Value = "5.0"
$ws.Range("C870").Style = "Normal"
$ws.Range("D870").Value = "A"
$ws.Range("F870").Value = "E4.Arg_A"
$ws.Range("G870").Value = "Argumenta las decisiones de diseño de manera sólida y coherente, basándose en fundamentos disciplinares y un profundo conocimiento del usuario. Las decisiones están bien justificadas y son consistentes."

